# Update countries & provincias Spain
# - Swap rank of Costa Rica / Singapur (Costa Rica now ahead, row 55 vs 56)
# - Swap rank of Cuba / Republica de Africa Central (Cuba now ahead, row 121 vs 122)
# - Refresh case numbers for several countries
# - Update the "datos actualizados" timestamp footer

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footer timestamp (row 1, col A) ---
$ws.Range("A1").Value = "Datos actualizados a 15 de Septiembre de 2020 a las 21:54"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 6772343
$ws.Range("C4").Value = 23054
$ws.Range("D4").Value = 4047571
$ws.Range("E4").Value = 2524977
$ws.Range("G4").Value = 795
$ws.Range("H4").Value = 199795

# --- Row 5: India ---
$ws.Range("B5").Value = 5018034
$ws.Range("C5").Value = 91120
$ws.Range("D5").Value = 3939111
$ws.Range("E5").Value = 996835

# --- Row 25: Alemania ---
$ws.Range("B25").Value = 264837
$ws.Range("C25").Value = 1616
$ws.Range("E25").Value = 17842
$ws.Range("G25").Value = 9
$ws.Range("H25").Value = 9445

# --- Rows 55-56: Costa Rica overtakes Singapur ---
$ws.Range("A55").Value = "Costa Rica"
$ws.Range("B55").Value = 58137
$ws.Range("C55").Value = 776
$ws.Range("D55").Value = 21536
$ws.Range("E55").Value = 35968
$ws.Range("G55").Value = 12
$ws.Range("H55").Value = 633

$ws.Range("A56").Value = "Singapur"
$ws.Range("B56").Value = 57488
$ws.Range("C56").Value = 34
$ws.Range("D56").Value = 56884
$ws.Range("E56").Value = 577
$ws.Range("H56").Value = 27

# --- Row 107: Mauritania ---
$ws.Range("B107").Value = 7319
$ws.Range("C107").Value = 24
$ws.Range("D107").Value = 6839
$ws.Range("E107").Value = 319

# --- Rows 121-122: Cuba overtakes Republica de Africa Central ---
$ws.Range("A121").Value = "Cuba"
$ws.Range("B121").Value = 4803
$ws.Range("C121").Value = 77
$ws.Range("D121").Value = 4119
$ws.Range("E121").Value = 576
$ws.Range("H121").Value = 108

$ws.Range("A122").Value = "Republica de Africa Central"
$ws.Range("B122").Value = 4772
$ws.Range("D122").Value = 1828
$ws.Range("E122").Value = 2882
$ws.Range("H122").Value = 62

# --- Row 128: Siria ---
$ws.Range("B128").Value = 3614
$ws.Range("C128").Value = 38
$ws.Range("D128").Value = 871
$ws.Range("E128").Value = 2583
$ws.Range("G128").Value = 3
$ws.Range("H128").Value = 160
